$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row ("Description" / "Price" on row 4) together with
# the three data rows (5-7) so nothing stray is left behind, then rebuild the
# three data rows as rows 1-3 (matching the target layout).
$ws.Rows("4:7").Delete()

# --- Row 1: barcode value "120000" | "T shirt" | "350000" ---
$ws.Range("A1").Value = "120000"
$ws.Range("A1").Font.Size = 16
$ws.Range("B1").Value = "T shirt"
$ws.Range("C1").Value = "350000"

# --- Row 2: "fhsjk93749" | "jeans" | 600000 ---
$ws.Range("A2").Value = "fhsjk93749"
$ws.Range("A2").Font.Size = 16
$ws.Range("B2").Value = "jeans"
# Column C is formatted as Text ("@"); round-trip through the "Normal" style
# so the number is actually stored as a number (matching the source
# workbook), not as text, while still ending up back on the Text format.
$ws.Range("C2").Style = "Normal"
$ws.Range("C2").Value = 600000
$ws.Range("C2").NumberFormat = "@"

# --- Row 3: 123456789 | "foular" | 120000 ---
$ws.Range("A3").Style = "Normal"
$ws.Range("A3").Value = 123456789
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Font.Size = 16

$ws.Range("B3").Value = "foular"

$ws.Range("C3").Style = "Normal"
$ws.Range("C3").Value = 120000
$ws.Range("C3").NumberFormat = "@"

# Match the row heights used by the other data rows (21pt).
$ws.Rows("1:3").RowHeight = 21

$ws.Range("C4").Select()
